$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated symbol list values (Price column D, Volume(1h) column E).
# Cells are stored as text in the source sheet, so force the Text number
# format before assigning to keep the values as literal strings rather
# than having Excel auto-convert numeric-looking / percent-looking text
# into actual numbers.
$changes = @{
    "D2" = "328.43"
    "E2" = "6.64%"
    "D3" = "39.97"
    "E3" = "6.25%"
    "E4" = "2.08%"
    "D5" = "0.08093"
    "E5" = "3.29%"
    "D6" = "4.533"
    "E6" = "2.25%"
    "E7" = "4.98%"
    "D8" = "1.933"
    "E8" = "1.76%"
    "E9" = "-1.15%"
    "D10" = "0.9366"
    "E10" = "0.44%"
    "E11" = "23.27%"
    "D12" = "0.1962"
    "E12" = "1.85%"
    "D13" = "0.09122"
    "E13" = "1.77%"
    "D14" = "0.03488"
    "E14" = "5.38%"
    "D15" = "0.09582"
    "E15" = "0.00%"
    "D16" = "0.001360"
    "E16" = "-1.85%"
    "D17" = "0.006209"
    "E17" = "7.34%"
    "E18" = "-6.88%"
    "D19" = "0.3509"
    "E19" = "3.58%"
    "D20" = "7.018"
    "E20" = "12.18%"
    "D21" = "0.1324"
    "E21" = "3.70%"
    "D22" = "0.2571"
    "E22" = "-0.51%"
    "D23" = "0.04432"
    "E23" = "1.12%"
    "E24" = "-0.70%"
    "D25" = "0.004274"
    "E25" = "-6.30%"
    "D26" = "0.0001191"
    "E26" = "-0.76%"
    "D27" = "0.0003994"
    "E27" = "0.06%"
    "D39" = "0.02501"
    "E39" = "14.06%"
    "D40" = "0.05179"
    "E40" = "3.02%"
    "D41" = "0.007698"
    "E41" = "3.46%"
    "D42" = "0.1423"
    "E42" = "5.55%"
    "D43" = "0.009071"
    "E43" = "3.95%"
    "E44" = "2.92%"
    "D45" = "0.01133"
    "E45" = "42.24%"
    "D46" = "0.00006649"
    "E46" = "1.91%"
    "E47" = "0.07%"
    "E48" = "148.16%"
    "D49" = "0.003346"
    "E49" = "17.00%"
    "E50" = "0.07%"
    "D51" = "0.0002002"
    "E51" = "0.07%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
